# Update TPM-derived NATMI LR-pair metrics for Thbs1-Sdc1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=7; Value=2.727484333333333},
    @{Row=2; Col=8; Value=8.182453},
    @{Row=2; Col=9; Value=0.03096049453772388},
    @{Row=2; Col=10; Value=0.03096049453772388},
    @{Row=2; Col=13; Value=0.2799683333333333},
    @{Row=2; Col=14; Value=0.839905},
    @{Row=2; Col=15; Value=0.0294305463214559},
    @{Row=2; Col=16; Value=0.0294305463214559},
    @{Row=2; Col=17; Value=0.7636092429961111},
    @{Row=2; Col=18; Value=6.872483186965001},
    @{Row=2; Col=19; Value=0.000911184268627665},
    @{Row=2; Col=20; Value=0.0009111842686276652},
    @{Row=3; Col=7; Value=2.727484333333333},
    @{Row=3; Col=8; Value=8.182453},
    @{Row=3; Col=9; Value=0.03096049453772388},
    @{Row=3; Col=10; Value=0.03096049453772388},
    @{Row=3; Col=15; Value=0.2486942046732164},
    @{Row=3; Col=16; Value=0.2486942046732163},
    @{Row=3; Col=17; Value=6.452656070118111},
    @{Row=3; Col=18; Value=58.073904631063},
    @{Row=3; Col=19; Value=0.0076996955653487},
    @{Row=3; Col=20; Value=0.0076996955653487},
    @{Row=4; Col=7; Value=2.727484333333333},
    @{Row=4; Col=8; Value=8.182453},
    @{Row=4; Col=9; Value=0.03096049453772388},
    @{Row=4; Col=10; Value=0.03096049453772388},
    @{Row=4; Col=13; Value=6.86709},
    @{Row=4; Col=14; Value=20.60127},
    @{Row=4; Col=15; Value=0.7218752490053277},
    @{Row=4; Col=16; Value=0.7218752490053277},
    @{Row=4; Col=17; Value=18.72988039059},
    @{Row=4; Col=18; Value=168.56892351531},
    @{Row=4; Col=19; Value=0.02234961470374752},
    @{Row=4; Col=20; Value=0.02234961470374752},
    @{Row=5; Col=9; Value=0.5986009007423507},
    @{Row=5; Col=10; Value=0.5986009007423507},
    @{Row=5; Col=13; Value=0.2799683333333333},
    @{Row=5; Col=14; Value=0.839905},
    @{Row=5; Col=15; Value=0.0294305463214559},
    @{Row=5; Col=16; Value=0.0294305463214559},
    @{Row=5; Col=17; Value=14.76388499271889},
    @{Row=5; Col=18; Value=132.87496493447},
    @{Row=5; Col=19; Value=0.01761715153736298},
    @{Row=5; Col=20; Value=0.01761715153736298},
    @{Row=6; Col=9; Value=0.5986009007423507},
    @{Row=6; Col=10; Value=0.5986009007423507},
    @{Row=6; Col=15; Value=0.2486942046732164},
    @{Row=6; Col=16; Value=0.2486942046732163},
    @{Row=6; Col=19; Value=0.1488685749267898},
    @{Row=6; Col=20; Value=0.1488685749267898},
    @{Row=7; Col=9; Value=0.5986009007423507},
    @{Row=7; Col=10; Value=0.5986009007423507},
    @{Row=7; Col=13; Value=6.86709},
    @{Row=7; Col=14; Value=20.60127},
    @{Row=7; Col=15; Value=0.7218752490053277},
    @{Row=7; Col=16; Value=0.7218752490053277},
    @{Row=7; Col=17; Value=362.12998015722},
    @{Row=7; Col=18; Value=3259.16982141498},
    @{Row=7; Col=19; Value=0.4321151742781979},
    @{Row=7; Col=20; Value=0.4321151742781979},
    @{Row=8; Col=7; Value=32.63402300000001},
    @{Row=8; Col=8; Value=97.90206900000001},
    @{Row=8; Col=9; Value=0.3704386047199253},
    @{Row=8; Col=10; Value=0.3704386047199253},
    @{Row=8; Col=13; Value=0.2799683333333333},
    @{Row=8; Col=14; Value=0.839905},
    @{Row=8; Col=15; Value=0.0294305463214559},
    @{Row=8; Col=16; Value=0.0294305463214559},
    @{Row=8; Col=17; Value=9.136493029271668},
    @{Row=8; Col=18; Value=82.22843726344502},
    @{Row=8; Col=19; Value=0.01090221051546525},
    @{Row=8; Col=20; Value=0.01090221051546525},
    @{Row=9; Col=7; Value=32.63402300000001},
    @{Row=9; Col=8; Value=97.90206900000001},
    @{Row=9; Col=9; Value=0.3704386047199253},
    @{Row=9; Col=10; Value=0.3704386047199253},
    @{Row=9; Col=15; Value=0.2486942046732164},
    @{Row=9; Col=16; Value=0.2486942046732163},
    @{Row=9; Col=17; Value=77.20525615117768},
    @{Row=9; Col=18; Value=694.8473053605991},
    @{Row=9; Col=19; Value=0.0921259341810778},
    @{Row=9; Col=20; Value=0.09212593418107778},
    @{Row=10; Col=7; Value=32.63402300000001},
    @{Row=10; Col=8; Value=97.90206900000001},
    @{Row=10; Col=9; Value=0.3704386047199253},
    @{Row=10; Col=10; Value=0.3704386047199253},
    @{Row=10; Col=13; Value=6.86709},
    @{Row=10; Col=14; Value=20.60127},
    @{Row=10; Col=15; Value=0.7218752490053277},
    @{Row=10; Col=16; Value=0.7218752490053277},
    @{Row=10; Col=17; Value=224.1007730030701},
    @{Row=10; Col=18; Value=2016.90695702763},
    @{Row=10; Col=19; Value=0.2674104600233823},
    @{Row=10; Col=20; Value=0.2674104600233823},
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
